$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text-valued updates (non-numeric-looking strings) - plain Value assignment keeps inline/shared string type
$ws.Range("D2").Value = "31.416.66"
$ws.Range("E2").Value = "  +3.37%  "
$ws.Range("D3").Value = "2.001.70"
$ws.Range("E3").Value = "  +6.92%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("E5").Value = "  +72.33%  "
$ws.Range("E6").Value = "  +4.86%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +24.02%  "
$ws.Range("E9").Value = "  +18.16%  "
$ws.Range("E10").Value = "  +8.26%  "
$ws.Range("E11").Value = "  +17.20%  "
$ws.Range("E12").Value = "  +4.65%  "
$ws.Range("E13").Value = "  +5.72%  "
$ws.Range("D14").Value = "1.993.85"
$ws.Range("E14").Value = "  +6.51%  "
$ws.Range("E15").Value = "  +7.55%  "
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "31.392.24"
$ws.Range("E17").Value = "  +3.33%  "
$ws.Range("E18").Value = "  +8.02%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("E19").Value = "  +11.01%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("E20").Value = "  +6.22%  "
$ws.Range("D21").Value = "2.258.87"
$ws.Range("E21").Value = "  +7.07%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("E24").Value = "  +13.47%  "
$ws.Range("E25").Value = "  +9.02%  "
$ws.Range("E26").Value = "  +60.76%  "
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("E28").Value = "  +7.31%  "
$ws.Range("E29").Value = "  +20.79%  "
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("E30").Value = "  +10.28%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("E31").Value = "  +6.23%  "
$ws.Range("E32").Value = "  +3.19%  "
$ws.Range("E33").Value = "  +6.12%  "
$ws.Range("E34").Value = "  +8.27%  "
$ws.Range("E35").Value = "  +8.91%  "
$ws.Range("E36").Value = "  +10.73%  "
$ws.Range("E37").Value = "  +1.59%  "
$ws.Range("E38").Value = "  +6.41%  "
$ws.Range("E39").Value = "  +3.66%  "
$ws.Range("E40").Value = "  +6.54%  "
$ws.Range("E41").Value = "  +13.06%  "
$ws.Range("E42").Value = "  +4.29%  "
$ws.Range("E43").Value = "  +10.54%  "
$ws.Range("E44").Value = "  +3.85%  "
$ws.Range("E45").Value = "  +3.77%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("E47").Value = "  +3.43%  "
$ws.Range("E48").Value = "  +8.23%  "
$ws.Range("E49").Value = "  +12.07%  "
$ws.Range("E50").Value = "  +4.85%  "
$ws.Range("E51").Value = "  +14.68%  "

# Numeric-looking string updates - must force text storage so Excel does not coerce them to numbers
# (matches the source data where these cells are stored as text, e.g. "0.9982")
$numericTextCells = @("D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D15","D16","D18","D19","D20","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).NumberFormat = "@"
}
$ws.Range("D4").Value = "0.9982"
$ws.Range("D5").Value = "0.8117"
$ws.Range("D6").Value = "256.01"
$ws.Range("D7").Value = "0.9982"
$ws.Range("D8").Value = "0.3580"
$ws.Range("D9").Value = "25.96"
$ws.Range("D10").Value = "0.07044"
$ws.Range("D11").Value = "0.8485"
$ws.Range("D12").Value = "0.08153"
$ws.Range("D13").Value = "101.75"
$ws.Range("D15").Value = "5.535"
$ws.Range("D16").Value = "274.12"
$ws.Range("D18").Value = "14.07"
$ws.Range("D19").Value = "5.853"
$ws.Range("D20").Value = "0.000007975"
$ws.Range("D22").Value = "0.9978"
$ws.Range("D23").Value = "0.9978"
$ws.Range("D24").Value = "7.090"
$ws.Range("D25").Value = "9.888"
$ws.Range("D26").Value = "0.1546"
$ws.Range("D27").Value = "164.71"
$ws.Range("D28").Value = "20.17"
$ws.Range("D29").Value = "2.275"
$ws.Range("D30").Value = "4.671"
$ws.Range("D31").Value = "1.580"
$ws.Range("D32").Value = "1.364"
$ws.Range("D33").Value = "4.377"
$ws.Range("D34").Value = "0.05234"
$ws.Range("D35").Value = "1.222"
$ws.Range("D36").Value = "0.7657"
$ws.Range("D37").Value = "2.758"
$ws.Range("D38").Value = "0.02021"
$ws.Range("D40").Value = "6.666"
$ws.Range("D41").Value = "0.4784"
$ws.Range("D42").Value = "79.03"
$ws.Range("D43").Value = "2.140"
$ws.Range("D44").Value = "0.8617"
$ws.Range("D45").Value = "104.75"
$ws.Range("D46").Value = "0.9988"
$ws.Range("D47").Value = "10.03"
$ws.Range("D48").Value = "7.562"
$ws.Range("D49").Value = "0.4414"
$ws.Range("D50").Value = "37.00"
$ws.Range("D51").Value = "0.1215"
foreach ($ref in $numericTextCells) {
    $ws.Range($ref).Style = "Normal"
}
